$wb = $excel.ActiveWorkbook

# ---- Sheet "Preços" (price sheet) : just a selection/cursor change ----
$ws1 = $wb.Worksheets.Item("Preços")
$ws1.Activate()
$ws1.Range("F9").Select()

# ---- Sheet "Planilha1" (data sheet) ----
$ws2 = $wb.Worksheets.Item("Planilha1")
$ws2.Activate()

# Data correction: Chuchu (row 11) type changed from "Verdura" to "Legume"
$ws2.Range("C11").Value = "Legume"

# New row 14: Salsinha / Verdura
$ws2.Range("B14").Value = "Salsinha"
$ws2.Range("C14").Value = "Verdura"
$ws2.Range("D14").Value = 8
$ws2.Range("E14").Value = 24
$ws2.Range("F14").Value = 10
$ws2.Range("G14").Value = 22
$ws2.Range("H14").Value = 75

# New row 15: Morango / Fruta
$ws2.Range("B15").Value = "Morango"
$ws2.Range("C15").Value = "Fruta"
$ws2.Range("D15").Value = 9
$ws2.Range("E15").Value = 30
$ws2.Range("F15").Value = 13
$ws2.Range("G15").Value = 26

# New row 16: Pimenta-vermelha (no Tipo)
$ws2.Range("B16").Value = "Pimenta-vermelha"
$ws2.Range("D16").Value = 18
$ws2.Range("E16").Value = 35
$ws2.Range("F16").Value = 20
$ws2.Range("G16").Value = 30

# Column B grew wider to fit the new "Pimenta-vermelha" label
$ws2.Columns.Item(2).ColumnWidth = 17

# Selection state on Planilha1 (this is the last-active sheet, matching the
# workbook's activeTab/tabSelected state)
$ws2.Range("G19").Select()
